$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.413.44'
$ws.Range('E2').Value = '  -0.78%  '
$ws.Range('D3').Value = '3.452.43'
$ws.Range('E3').Value = '  -1.54%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '609.63'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.23%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '167.57'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.20%  '
$ws.Range('D7').Value = '3.444.69'
$ws.Range('E7').Value = '  -1.63%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.595'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.19%  '
$ws.Range('E9').Value = '  +0.12%  '
$ws.Range('E10').Value = '  -0.03%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '7.04'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.41%  '
$ws.Range('E12').Value = '  -2.97%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '44.32'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.42%  '
$ws.Range('E14').Value = '  -1.80%  '
$ws.Range('D15').Value = '4.014.30'
$ws.Range('E15').Value = '  -1.30%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '8.17'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.29%  '
$ws.Range('B17').Value = 'WrappedBTC'
$ws.Range('C17').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D17').Value = '69.543.99'
$ws.Range('E17').Value = '  -0.72%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D18').Value = '3.449.04'
$ws.Range('E18').Value = '  -1.54%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '576.72'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.87%  '
$ws.Range('E20').Value = '  +1.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.846'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.82%  '
$ws.Range('E23').Value = '  -1.79%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '95.93'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.44%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '15.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.69%  '
$ws.Range('E26').Value = '  -2.23%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.12%  '
$ws.Range('E28').Value = '  -5.07%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '32.79'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -3.17%  '
$ws.Range('E30').Value = '  -3.86%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.89%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.81'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -6.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.24'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.72%  '
$ws.Range('E34').Value = '  -5.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '578.58'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -18.58%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '10.54'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.30%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0475'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.85%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0955'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.21%  '
$ws.Range('E39').Value = '  +0.26%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '56.05'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.60%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.140'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.14'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -11.32%  '
$ws.Range('D43').Value = '3.249.72'
$ws.Range('E43').Value = '  -2.35%  '
$ws.Range('D44').Value = '0.0₃0689'
$ws.Range('E44').Value = '  +0.21%  '
$ws.Range('E45').Value = '  -5.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '31.05'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.79'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -4.56%  '
$ws.Range('E48').Value = '  -5.99%  '
$ws.Range('E49').Value = '  -2.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '133.89'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.43%  '
